# Loan RBI, Variable Instalments
# Insert a new (blank) "Variable Instalment" column into the Repayment
# schedule sheet, shifting the existing "Late" / "Outstanding" columns one
# place to the right, and leave that sheet active/selected as the user left
# it after making the edit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"): this pushes the
# existing N/O/P (Late / heading / Outstanding) columns to O/P/Q.
$ws.Columns("N").Insert() | Out-Null

# The inserted column keeps the same width as its neighbour (column M).
$ws.Columns("N").ColumnWidth = 10.17

# The "Repayment schedule" tab is the one left selected/active, with the
# cursor resting on J16.
$ws.Activate() | Out-Null
$ws.Range("J16").Select() | Out-Null
